# Auto-update draw results: append the 2025-11-09 Pick 3 draw as a new
# row 54 (dimension grows from A1:E53 to A1:E54).
#
# Columns A (date "2025-11-09") and C (phase "251109") look numeric/date-like,
# so a plain .Value assignment would coerce them to a date serial / number.
# The source data stores every column as literal text, so we force text
# entry the same way a user would in Excel - a leading apostrophe - and
# then reset the cell Style back to "Normal" so no left-over
# quote-prefix/number-format styling is attached to the cell (matching the
# unstyled text cells used throughout the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "'2025-11-09"
$ws.Range("B54").Value = "Pick 3"
$ws.Range("C54").Value = "'251109"
$ws.Range("D54").Value = "1-7-6"
$ws.Range("E54").Value = "2025-11-09T21:36:18.004+04:00"

$ws.Range("A54").Style = "Normal"
$ws.Range("C54").Style = "Normal"
